# Insert two new rows right before the current row 63 ("Lapins / Región de Ñuble"
# 12-Jan-2021 entries). This pushes the existing rows 63..84 down to 65..86,
# preserving their values/styles, and grows the used range to A1:T86.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("63:64").Insert()

# --- New row 63: Lapins / Primera, 29-Dec-2021, Región de O'Higgins ---
$ws.Cells.Item(63, 1).Value = 11
$ws.Cells.Item(63, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(63, 3).Value = "Bíobío"
$ws.Cells.Item(63, 4).Value = 44559
$ws.Cells.Item(63, 5).Value = 8
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100103
$ws.Cells.Item(63, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(63, 9).Value = 100103001
$ws.Cells.Item(63, 10).Value = "Cereza"
$ws.Cells.Item(63, 11).Value = "Lapins"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 100
$ws.Cells.Item(63, 14).Value = 4500
$ws.Cells.Item(63, 15).Value = 5000
$ws.Cells.Item(63, 16).Value = 4750
$ws.Cells.Item(63, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(63, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(63, 19).Value = 475
$ws.Cells.Item(63, 20).Value = 10

# --- New row 64: Lapins / Segunda, 29-Dec-2021, Región de O'Higgins ---
$ws.Cells.Item(64, 1).Value = 11
$ws.Cells.Item(64, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(64, 3).Value = "Bíobío"
$ws.Cells.Item(64, 4).Value = 44559
$ws.Cells.Item(64, 5).Value = 8
$ws.Cells.Item(64, 6).Value = "Fruta"
$ws.Cells.Item(64, 7).Value = 100103
$ws.Cells.Item(64, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(64, 9).Value = 100103001
$ws.Cells.Item(64, 10).Value = "Cereza"
$ws.Cells.Item(64, 11).Value = "Lapins"
$ws.Cells.Item(64, 12).Value = "Segunda"
$ws.Cells.Item(64, 13).Value = 50
$ws.Cells.Item(64, 14).Value = 4000
$ws.Cells.Item(64, 15).Value = 4000
$ws.Cells.Item(64, 16).Value = 4000
$ws.Cells.Item(64, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(64, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(64, 19).Value = 400
$ws.Cells.Item(64, 20).Value = 10

# Keep the date format that row 2..84 already use (style index 2 -> date numFmt)
$ws.Range("D63:D64").NumberFormat = $ws.Range("D65").NumberFormat
